# Actualización automática: registra la venta de noviembre del cliente
# "ROCA REYNA PAUL DAVID" (asesor LOZANO MOLINA TITO) y refleja el
# impacto en los totales / contadores de las hojas relacionadas.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Hoja "VENTAS POR GRUPO": detalle de venta por grupo de producto
# ---------------------------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Fila 24 = ROCA REYNA PAUL DAVID
$wsGrupo.Range("C24").Value = 259.2    # 240X120 PORCELANATO
$wsGrupo.Range("I24").Value = 298.8    # LAVABOS
$wsGrupo.Range("L24").Value = 665.16   # PIEDRA SINTERIZADA
$wsGrupo.Range("N24").Value = 231.88   # PUERTAS DE SEGURIDAD

# Fila 32 = contador "X de 30" clientes que compraron cada grupo
$wsGrupo.Range("C32").Value = "1 de 30"
$wsGrupo.Range("I32").Value = "1 de 30"
$wsGrupo.Range("L32").Value = "3 de 30"
$wsGrupo.Range("N32").Value = "1 de 30"

# ---------------------------------------------------------------
# Hoja "VENTA MENSUAL": total de venta de noviembre del cliente
# ---------------------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

$wsMensual.Range("F24").Value = 1455.04   # noviembre - ROCA REYNA PAUL DAVID
$wsMensual.Range("F32").Value = 11263.58  # noviembre - total general
